$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 29-31 are being reordered: InjectiveProtocol, Cosmos, Kaspa (were Cosmos, Kaspa, InjectiveProtocol)
# plus general price/volume updates across the sheet.

# Row 2
$ws.Range("D2").Value = '49.990.31'
$ws.Range("E2").Value = '  +4.08%  '

# Row 3
$ws.Range("D3").Value = '2.638.54'
$ws.Range("E3").Value = '  +5.54%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.46'
$ws.Range("E5").Value = '  +2.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.62'
$ws.Range("E6").Value = '  +3.06%  '

# Row 7
$ws.Range("E7").Value = '  +1.74%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  +4.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.76'
$ws.Range("E10").Value = '  +3.32%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.80'
$ws.Range("E11").Value = '  +3.17%  '

# Row 12
$ws.Range("E12").Value = '  +1.39%  '

# Row 13
$ws.Range("E13").Value = '  +0.92%  '

# Row 14
$ws.Range("E14").Value = '  +2.93%  '

# Row 15
$ws.Range("D15").Value = '3.052.88'
$ws.Range("E15").Value = '  +5.61%  '

# Row 16
$ws.Range("D16").Value = '2.621.21'
$ws.Range("E16").Value = '  +5.09%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.883'
$ws.Range("E17").Value = '  +5.53%  '

# Row 18
$ws.Range("D18").Value = '49.929.81'
$ws.Range("E18").Value = '  +4.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.08'
$ws.Range("E19").Value = '  +11.98%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.39'
$ws.Range("E20").Value = '  +3.72%  '

# Row 21
$ws.Range("E21").Value = '  +1.76%  '

# Row 22
$ws.Range("E22").Value = '  +2.80%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.08'
$ws.Range("E23").Value = '  +2.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '281.27'
$ws.Range("E24").Value = '  +1.42%  '

# Row 25
$ws.Range("E25").Value = '  +2.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.68'
$ws.Range("E26").Value = '  +4.30%  '

# Row 27
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("E28").Value = '  +7.47%  '

# Row 29
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.64'
$ws.Range("E29").Value = '  +4.60%  '

# Row 30
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.95'
$ws.Range("E30").Value = '  +2.38%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.144'
$ws.Range("E31").Value = '  +4.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.82'
$ws.Range("E32").Value = '  +0.78%  '

# Row 33
$ws.Range("E33").Value = '  +1.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.44'
$ws.Range("E34").Value = '  +2.92%  '

# Row 35
$ws.Range("E35").Value = '  -0.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0798'
$ws.Range("E36").Value = '  +2.36%  '

# Row 37
$ws.Range("E37").Value = '  +6.58%  '

# Row 38
$ws.Range("E38").Value = '  +3.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("E39").Value = '  +8.09%  '

# Row 40
$ws.Range("E40").Value = '  +1.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '123.70'
$ws.Range("E41").Value = '  +2.58%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.55'
$ws.Range("E42").Value = '  +6.44%  '

# Row 43
$ws.Range("E43").Value = '  +0.39%  '

# Row 44
$ws.Range("E44").Value = '  +4.86%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +6.74%  '

# Row 46
$ws.Range("D46").Value = '2.064.56'
$ws.Range("E46").Value = '  +2.93%  '

# Row 47
$ws.Range("E47").Value = '  +16.00%  '

# Row 48
$ws.Range("E48").Value = '  +8.50%  '

# Row 49
$ws.Range("E49").Value = '  +1.25%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.40'
$ws.Range("E50").Value = '  +4.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.76'
$ws.Range("E51").Value = '  +2.07%  '
